$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 0.3879434203120516
$ws.Range("B2").Value = 0.2344599303095528
$ws.Range("C2").Value = 0.3816025998094102

$ws.Range("A3").Value = 0.3613924994661383
$ws.Range("B3").Value = 0.398713354387823
$ws.Range("C3").Value = 0.2393950352037498

$ws.Range("A4").Value = 0.3875483107160846
$ws.Range("B4").Value = 0.3923488752478406
$ws.Range("C4").Value = 0.2253523643097708

$ws.Range("A5").Value = 0.3895877313122079
$ws.Range("B5").Value = 0.390511918804573
$ws.Range("C5").Value = 0.2198814110889392

$ws.Range("A6").Value = 0.3958558681312614
$ws.Range("B6").Value = 0.2221436593682502
$ws.Range("C6").Value = 0.3820490584326038
